# Updated cryptos list on Mon Jun 17 18:29:18 UTC 2024 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto ranking table on the active sheet.
#
# All Price/Volume cells in the source sheet are stored as TEXT (inline
# strings), not numbers - e.g. thousands separators use "." (like
# "66.643.79") and percentages keep two leading/trailing spaces
# (e.g. "  +0.08%  "). A bare numeric-looking assignment such as
# $ws.Range("D4").Value = "1.00" would be auto-coerced by Excel into the
# number 1, dropping the trailing zero and changing the cell type - so any
# new Price value that looks like a plain number is written with a leading
# apostrophe to force it to stay text, exactly as typing '1.00 into Excel
# would. Non-numeric-looking Price strings (those keeping the "."
# thousands separator, e.g. "67.000.17") and all Volume percentages are
# already safe to assign directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.000.17'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '3.564.00'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''608.73'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").Value = '''145.92'
$ws.Range("E6").Value = '  -1.55%  '
$ws.Range("D7").Value = '3.564.80'
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = '''0.512'
$ws.Range("E9").Value = '  +4.76%  '
$ws.Range("D10").Value = '''7.88'
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("D11").Value = '''0.133'
$ws.Range("E11").Value = '  -2.63%  '
$ws.Range("D12").Value = '''0.413'
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").Value = '4.172.23'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").Value = '''0.0000196'
$ws.Range("E14").Value = '  -6.29%  '
$ws.Range("D15").Value = '''28.99'
$ws.Range("E15").Value = '  -3.01%  '
$ws.Range("D16").Value = '3.562.32'
$ws.Range("D17").Value = '''0.118'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '66.759.67'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").Value = '''11.12'
$ws.Range("E19").Value = '  -2.77%  '
$ws.Range("D20").Value = '''6.21'
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").Value = '''14.72'
$ws.Range("E21").Value = '  -2.58%  '
$ws.Range("D22").Value = '''427.04'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '''0.599'
$ws.Range("E23").Value = '  -3.36%  '
$ws.Range("D24").Value = '''77.68'
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("D25").Value = '3.706.51'
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -3.55%  '
$ws.Range("D28").Value = '''8.04'
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("D29").Value = '''2.49'
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("D30").Value = '''9.05'
$ws.Range("E30").Value = '  -2.98%  '
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").Value = '3.572.48'
$ws.Range("E32").Value = '  -0.48%  '
$ws.Range("D33").Value = '''0.157'
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("D34").Value = '''24.46'
$ws.Range("E34").Value = '  -4.05%  '
$ws.Range("D36").Value = '''1.36'
$ws.Range("E36").Value = '  -6.77%  '
$ws.Range("D37").Value = '''7.69'
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("E38").Value = '  -3.21%  '
$ws.Range("D39").Value = '''177.55'
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").Value = '''5.33'
$ws.Range("E40").Value = '  -5.46%  '
$ws.Range("D41").Value = '''0.0828'
$ws.Range("E41").Value = '  -3.52%  '
$ws.Range("D42").Value = '''5.03'
$ws.Range("E42").Value = '  -3.92%  '
$ws.Range("D43").Value = '''0.867'
$ws.Range("E43").Value = '  -3.51%  '
$ws.Range("D44").Value = '''45.55'
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = '''1.79'
$ws.Range("E45").Value = '  -6.31%  '
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '''2.44'
$ws.Range("E47").Value = '  -4.98%  '
$ws.Range("D48").Value = '''23.76'
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("D49").Value = '''7.14'
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("E50").Value = '  -4.88%  '
$ws.Range("D51").Value = '''0.920'
$ws.Range("E51").Value = '  -3.26%  '
